$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("variables")

# --- New row of data: Samoa8 / 2,3 / Somewhat or very proud of being Samoan / Culture ---
# Clear any inherited styling on the row first so the new cells carry the default style
$ws.Range("A5:D5").Style = "Normal"

# Values are entered in this order so the shared-string table is built in the same
# order as the source workbook (title column first, then var/val/varname).
$ws.Range("D5").Value = "Culture"
$ws.Range("A5").Value = "Samoa8"
$ws.Range("B5").Value = "2,3"
$ws.Range("C5").Value = "Somewhat or very proud of being Samoan"

# --- Conditional formatting ---
# A2:A4 and A6 keep watching for ISFORMULA(); only their priority shifts down because a
# new, higher-priority rule is being inserted for the now-populated row 5. Updating the
# Priority property in place keeps each rule's existing style (dxf) association intact.
$ws.Range("A2:A4").FormatConditions.Item(1).Priority = 24
$ws.Range("A6").FormatConditions.Item(1).Priority = 7

# A5 no longer needs the "is this a formula" check (it now holds real data), so drop it...
$ws.Range("A5").FormatConditions.Delete()

# ...and replace it with validation rules appropriate for manually entered data:
# 1) flag the cell if it is left blank
$blankRule = $ws.Range("A5").FormatConditions.Add(2, 0, "LEN(TRIM(A5))=0")
$blankRule.StopIfTrue = $true
$blankRule.Priority = 2

# 2) flag the cell if its value duplicates another variable name already used
$dupRule = $ws.Range("A5").FormatConditions.AddUniqueValues()
$dupRule.DupeUnique = 1
$dupRule.Priority = 1
$dupRule.Font.Color = 0x0006009C
$dupRule.Interior.Color = 0x00CEC7FF

# --- Window / selection state ---
$ws.Range("F5").Select()

$wb.Save()
